$d = $word.ActiveDocument

# 1. Update the Heading1 title paragraph: split into title + br + URL
$d.Content.Find.Execute("Review 148: Knowledge Graph Prompting for Multi-Document Question Answering, 17.09.23", $true, $false, $false, $false, $false, $true, 1, $false, "Review 147: Generative Image Dynamics^lhttps://huggingface.co/papers/2309.07906", 2) | Out-Null

# 2. Update the bold "Paper:" line with the new arXiv link
$d.Content.Find.Execute("Paper: https://arxiv.org/abs/2308.11730v3", $true, $false, $false, $false, $false, $true, 1, $false, "Paper: https://arxiv.org/abs/2309.07906v3", 2) | Out-Null

# 3. Remove the paragraph with the old arXiv link + Hebrew credit line entirely
$p4 = $d.Paragraphs.Item(4)
$p4.Range.Delete()

# 4. Replace the body text of the first content paragraph (now index 5) with the new review text
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Find.Execute("אנחנו יודעים שמודלי שפה יודעים לענות על שאלות על גבי מסמכים אם בצורת zero shot או בצורת few shot אבל מה עם המקרה שבו צריך להסתכל על מסמכים שונים כדי לקבל תשובה? או להשוות בין שני מסמכים בעלי מבנה שונה (טבלה וטקסט חופשי)? ספוילר, לא טוב. היום ב-#shorthebrewpapereviews קבלו את הסקירה הבאה:", $true, $false, $false, $false, $false, $true, 1, $false, "נתחיל מזה שזה מאמר מאוד יפה רק בגלל זה מופיע בו התמרת פוריה (עבדתי בתחום עיבוד אותות כמה שנים ויש לי זכרונות חמים ממנו). שנית, המאמר עוסק במודלי דיפוזיה וביצירת וידאו עם תמונה שאני מאוד אוהב. אוקיי.^l^lאז היום ב-#shorthebrewpapereviews סוקרים מאמר ש״מוסיף״ תנועה לאובייקט בתמונה, כלומר מגנרט וידאו קצר המראה אובייקט בדינמיקה( למשל פרח ברוח).  יש כמה גישות למידול של תנועה בוידאו למשל שדה תנועה שמסמן תזוזה של כל פיקסל בתמונה. כלומר וידאו המורכב מ-T פריימים ניתן לייצוג על ידי T x H x W מספרים לתמונה בגודל HxW. ^l^lכלומר יש לנו  HxW מערכים באורך T מספרים המתארים את המסלול של כל פיקסל בתמונה. אז הנה באה הקטע של התמרת פוריה שאני כה אוהב. במקום לחזות את המערך הזה בואו נעשה זאת עם התמרת הפוריה שלו. כמובן שבמקרה שלנו ההתמרה דו מימדית כי כל פיקסל יכול לנוע בשני הכיוונים (למעשה 4 מספרים ממשים פר פיקסל כי התמרת פוריה היא מרוכבת). ^l^lאוקיי, אז מאמנים מודל דיפוזיה שמטרתו לגנרט התמרת פוריה של תנועתו של כל פיקסל. הבעיה אבל שקצת קשה לאמן רשת שחוזה וקטור ארוך (של תדרים) לכל פיקסל בתמונה (לטענת המאמר). אז הם שמים לב שתנועה של כל פיקסל מתוארת בעיקר על ידי תדרים נמוכים (תנועה איטית) כאשר התדרים הגבוהים מקבלים ערכים נמוכים. אז המחברים מחליטים לקחת רק 16 התדרים הנמוכים ביותר לגנרוט (יש איזה נרמול של התדרים כדי למנוע ערכים גבוהים מדי של תדרים מסוימים). ^l^lטוב, איך מאמנים מודל דיפוזיה. לוקחים תמונה בדינמיקה (סרטון קצר), מחשבים את שדה התמונה, מעבירים את זה דרך התמרת פוריה וזה מהווה הקלט למודל דיפוזיה שאותו מרעישים ואז מסירים את הרעש לאט לאט עם המודל. בנוסף מכניסים למודל את התמונה וגם תיאור טקסטואלי (אופציאונלי). ב-inference מעבירים את התמרת פוריה של שדה התנועה דרך התמרת פוריה הפוכה ומרנדרים את הוידאו. והדבר האחרון כמובן משתמשים במודל דיפוזיה לטנטי….", 2) | Out-Null

# 5. Delete the remaining old paragraphs (now indices 6 through 12) that held the rest of the old review text
$pStart = $d.Paragraphs.Item(6)
$pEnd = $d.Paragraphs.Item(12)
$delRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$delRange.Delete()

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
